$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value = "Jashanpreet Sidhu"

# Reusable text blocks (Method Inputs column F)
$f7  = "account_number = 2000`nclient_number = 2000`nbalance = 2000`ndate_created = (2015, 1, 1)`nmanagement_fee = 10"
$f8  = "account_number = 2000`nclient_number = 2000`nbalance = 2000`ndate_created = (2015, 1, 1)`nmanagement_fee = 'ten'"
$f9  = "account_number = 2000`nclient_number = 2000`nbalance = 2000`ndate_created = (2010, 1, 1)`nmanagement_fee = 10"
$f10 = "account_number = 2000`nclient_number = 2000`nbalance = 2000`ndate_created = (2015, 2, 16)`nmanagement_fee = 10"

# Reusable text blocks (Expected Result column G)
$g7  = "attributes are set "
$g8  = "management_fee set to 2.55"
$g9  = "service_charge set to base charge"
$g10 = "service_charge calculated on the basis of formula"
$g12 = "returns formatted string"

# Row 7 filled in first, completely
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = $f7
$ws.Range("G7").Value = $g7

# Then Method Inputs (F) filled down for rows 8-10
$ws.Range("F8").Value  = $f8
$ws.Range("F9").Value  = $f9
$ws.Range("F10").Value = $f10

# Then Expected Result (G) filled down for rows 8-10
$ws.Range("G8").Value  = $g8
$ws.Range("G9").Value  = $g9
$ws.Range("G10").Value = $g10

# Condition being Tested (E) filled down for rows 8-13
$ws.Range("E8").Value  = "None"
$ws.Range("E9").Value  = "None"
$ws.Range("E10").Value = "None"
$ws.Range("E11").Value = "None"

# Row 11 reuses row 7's input and row 10's expected result
$ws.Range("F11").Value = $f7
$ws.Range("G11").Value = $g10

# Row 12 reuses row 9's input, introduces a new expected result
$ws.Range("E12").Value = "None"
$ws.Range("F12").Value = $f9
$ws.Range("G12").Value = $g12

# Row 13 reuses row 7's input and row 12's expected result
$ws.Range("E13").Value = "None"
$ws.Range("F13").Value = $f7
$ws.Range("G13").Value = $g12

# Writing multi-line text above triggered row auto-fit; restore the
# original (manually set) row heights for rows 7-11 so they keep the
# explicit customHeight previously set on these rows.
$ws.Rows.Item(7).RowHeight  = 31.15
$ws.Rows.Item(8).RowHeight  = 31.15
$ws.Rows.Item(9).RowHeight  = 31.15
$ws.Rows.Item(10).RowHeight = 31.15
$ws.Rows.Item(11).RowHeight = 31.15

# Selection moved to G14 (matches the final cursor position recorded in the saved file)
$ws.Range("G14").Select()
